$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename headers to add "(AC)" suffix for the AC-based simulation tools
$ws.Range("B1").Value = "ESP(AC)"
$ws.Range("C1").Value = "BLAST(AC)"
$ws.Range("D1").Value = "DOE2(AC)"
$ws.Range("E1").Value = "SRES/SUN(AC)"
$ws.Range("F1").Value = "SERIRES(AC)"
$ws.Range("G1").Value = "S3PAS(AC)"
$ws.Range("H1").Value = "TRNSYS(AC)"
$ws.Range("I1").Value = "TASE(AC)"

# Minor value fixes in the distance_% column (M)
$ws.Range("M11").Value = 0.1
$ws.Range("M18").Value = -14.6
